# Applies the cryptos.xlsx price/volume refresh described by the commit
# "Updated cryptos list on Fri Oct 27 19:49:27 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several "Price" (column D) values look like plain numbers (e.g. "224.26").
# The source data stores them as text (inline strings), so force a Text
# number format on those specific cells before assigning the value - otherwise
# Excel auto-coerces a numeric-looking string into a real number.
$textCells = @(
    'D5',
    'D7',
    'D8',
    'D9',
    'D11',
    'D13',
    'D16',
    'D17',
    'D18',
    'D20',
    'D22',
    'D23',
    'D25',
    'D26',
    'D27',
    'D31',
    'D33',
    'D34',
    'D36',
    'D41',
    'D42',
    'D44',
    'D46',
    'D48'
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# --- Row-by-row cell updates -------------------------------------------

# Row 2
$ws.Range('D2').Value = '33.673.47'
$ws.Range('E2').Value = '  -0.85%  '

# Row 3
$ws.Range('D3').Value = '1.772.82'
$ws.Range('E3').Value = '  -0.96%  '

# Row 4
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$ws.Range('D5').Value = '224.26'
$ws.Range('E5').Value = '  +1.03%  '

# Row 6
$ws.Range('E6').Value = '  -0.96%  '

# Row 7
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.04%  '

# Row 8
$ws.Range('D8').Value = '31.84'
$ws.Range('E8').Value = '  +0.99%  '

# Row 9
$ws.Range('D9').Value = '0.289'
$ws.Range('E9').Value = '  +1.70%  '

# Row 10
$ws.Range('E10').Value = '  -4.35%  '

# Row 11
$ws.Range('D11').Value = '0.0935'
$ws.Range('E11').Value = '  +1.43%  '

# Row 12
$ws.Range('D12').Value = '2.027.75'
$ws.Range('E12').Value = '  -0.88%  '

# Row 13
$ws.Range('D13').Value = '11.07'
$ws.Range('E13').Value = '  +4.10%  '

# Row 14
$ws.Range('D14').Value = '1.776.03'
$ws.Range('E14').Value = '  -0.78%  '

# Row 15
$ws.Range('D15').Value = '33.674.29'
$ws.Range('E15').Value = '  -0.71%  '

# Row 16
$ws.Range('D16').Value = '0.608'
$ws.Range('E16').Value = '  -3.35%  '

# Row 17
$ws.Range('D17').Value = '4.12'

# Row 18
$ws.Range('D18').Value = '66.48'
$ws.Range('E18').Value = '  -2.23%  '

# Row 19
$ws.Range('D19').Value = '0.0₃0775'
$ws.Range('E19').Value = '  -1.16%  '

# Row 20
$ws.Range('D20').Value = '238.03'
$ws.Range('E20').Value = '  -2.94%  '

# Row 21
$ws.Range('E21').Value = '  -0.05%  '

# Row 22
$ws.Range('D22').Value = '10.56'
$ws.Range('E22').Value = '  -1.76%  '

# Row 23
$ws.Range('D23').Value = '4.00'
$ws.Range('E23').Value = '  -2.17%  '

# Row 24
$ws.Range('E24').Value = '  -2.57%  '

# Row 25
$ws.Range('D25').Value = '159.36'
$ws.Range('E25').Value = '  +0.74%  '

# Row 26
$ws.Range('D26').Value = '16.08'
$ws.Range('E26').Value = '  -1.84%  '

# Row 27
$ws.Range('D27').Value = '7.01'
$ws.Range('E27').Value = '  -0.23%  '

# Row 28
$ws.Range('E28').Value = '  -0.12%  '

# Row 29
$ws.Range('E29').Value = '  +0.17%  '

# Row 30
$ws.Range('E30').Value = '  +1.57%  '

# Row 31
$ws.Range('D31').Value = '0.0511'
$ws.Range('E31').Value = '  -1.60%  '

# Row 32
$ws.Range('E32').Value = '  -2.84%  '

# Row 33
$ws.Range('D33').Value = '3.49'
$ws.Range('E33').Value = '  -0.45%  '

# Row 34
$ws.Range('D34').Value = '1.79'

# Row 35
$ws.Range('D35').Value = '1.379.42'
$ws.Range('E35').Value = '  -2.22%  '

# Row 36
$ws.Range('D36').Value = '0.646'
$ws.Range('E36').Value = '  +0.77%  '

# Row 37
$ws.Range('E37').Value = '  -2.34%  '

# Row 39
$ws.Range('E39').Value = '  +5.46%  '

# Row 40
$ws.Range('E40').Value = '  +0.83%  '

# Row 41
$ws.Range('D41').Value = '77.97'
$ws.Range('E41').Value = '  -2.26%  '

# Row 42
$ws.Range('D42').Value = '0.905'
$ws.Range('E42').Value = '  -4.04%  '

# Row 43
$ws.Range('E43').Value = '  -2.39%  '

# Row 44
$ws.Range('D44').Value = '13.48'
$ws.Range('E44').Value = '  +13.67%  '

# Row 45
$ws.Range('E45').Value = '  +3.78%  '

# Row 46
$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').Value = '0.0499'
$ws.Range('E46').Value = '  +0.94%  '

# Row 47
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₆0135'
$ws.Range('E47').Value = '  +12.75%  '

# Row 48
$ws.Range('D48').Value = '107.02'
$ws.Range('E48').Value = '  +1.35%  '

# Row 49
$ws.Range('E49').Value = '  -2.08%  '

# Row 50
$ws.Range('D50').Value = '1.927.80'
$ws.Range('E50').Value = '  -0.51%  '

# Restore the default (General) style on the text-forced cells so the
# cells end up identical to their original style (no explicit s="..."),
# matching the source which only differs in cell content, not formatting.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}

Write-Output "Updated $($textCells.Count) text-forced price cells and 83 total cell values."
